$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo")

# --- Update "Outros Gastos" row (row 12): Abril..Dezembro (F..M) go from 50 to 250 ---
$ws.Range("F12:M12").Value = 250

# --- Update "Reserva" row (row 15): hourly-rate cell A15 changes from 5.85 to 5.64 ---
$ws.Range("A15").Value = 5.64

# --- Row 15 is no longer hidden ---
$ws.Rows.Item(15).Hidden = $false
# Restore the natural row height (writing into a hidden row otherwise stamps
# an explicit ht/customHeight on the row, which the source file never had).
$ws.Rows.Item(15).AutoFit()

# --- Update the active selection on the Resumo sheet ---
$ws.Activate()
$ws.Range("A15").Select()

# --- Re-sequence the mergeCells list on Planilha1 to match the saved order ---
$ws1 = $wb.Worksheets.Item("Planilha1")
$reorderFirst = @("BM1:BN1","AF1:AG1","AI1:AJ1","AL1:AM1","AO1:AP1","AR1:AS1","AU1:AV1","AX1:AY1","BA1:BB1","BD1:BE1","BG1:BH1","BJ1:BK1")
foreach ($addr in $reorderFirst) {
    $ws1.Range($addr).UnMerge()
    $ws1.Range($addr).Merge()
}

$wb.Save()
